$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - Enterprises density (per 1000 people)
$ws.Range("B13").Value = "'34.13"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'1.59"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'35.71"
$ws.Range("D13").Style = "Normal"

# Row 14 - Employment (% of total)
$ws.Range("B14").Value = "'44.78"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'36.46"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'81.24"
$ws.Range("D14").Style = "Normal"

# Row 16 - Enterprises (% of total)
$ws.Range("B16").Value = "'95.44"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'4.44"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'99.88"
$ws.Range("D16").Style = "Normal"
